# KPI Template v2 - add Base Measurement and Bay Count sheets,
# rename store_fk -> store_id on KPIS sheet.

$wb = $excel.ActiveWorkbook
$kpis = $wb.Worksheets.Item("KPIS")

# ---------------------------------------------------------------
# 1. Rename the "store_fk" field label to "store_id" everywhere it
#    appears on the KPIS sheet (shared text string).
# ---------------------------------------------------------------
[void]$kpis.Cells.Replace("store_fk", "store_id", -4163, 1, $false, $false, $false)

# Reference cells used to copy existing formatting from, so the new
# sheets look consistent with the rest of the workbook.
$shareOfShelf   = $wb.Worksheets.Item("Share of Shelf")
$headerStyleSrc = $kpis.Range("A1")         # bold / green header style
$plainStyleSrc  = $shareOfShelf.Range("F2") # plain font, no border
$borderStyleSrc = $kpis.Range("A12")        # plain font, thin border all round

# ---------------------------------------------------------------
# 2. Add the "Base Measurement" sheet (after "Hierarchy").
# ---------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$baseMeasurement = $wb.Worksheets.Add($null, $lastSheet)
$baseMeasurement.Name = "Base Measurement"

$baseMeasurement.Range("A1").Value2 = "KPI Name"
$baseMeasurement.Range("B1").Value2 = "Sum Col"

$baseMeasurement.Range("A2").Value2 = "Warm Base Measurement"
$baseMeasurement.Range("B2").Value2 = "net_len_ign_stack"

$baseMeasurement.Range("A3").Value2 = "Cold Room Base Measurement"
$baseMeasurement.Range("B3").Value2 = "net_len_ign_stack"

$baseMeasurement.Range("A4").Value2 = "Cooler Door Measurement"
$baseMeasurement.Range("B4").Value2 = "net_len_ign_stack"

# formatting: column A rows 2-4 get a thin border (like other KPI sheets)
$borderStyleSrc.Copy() | Out-Null
$baseMeasurement.Range("A2:A4").PasteSpecial(-4122) | Out-Null

# column B row 2 stays plain, rows 3-4 use the "applyFont" plain style
$plainStyleSrc.Copy() | Out-Null
$baseMeasurement.Range("B3:B4").PasteSpecial(-4122) | Out-Null

$baseMeasurement.Columns.Item(1).ColumnWidth = 26.5026990553307
$baseMeasurement.Columns.Item(2).ColumnWidth = 14.8184885290149

$baseMeasurement.Rows.Item(1).RowHeight = 12.8
$baseMeasurement.Rows.Item(2).RowHeight = 13.8
$baseMeasurement.Rows.Item(3).RowHeight = 13.8
$baseMeasurement.Rows.Item(4).RowHeight = 13.8

[void]$baseMeasurement.Range("A1").Select()

# ---------------------------------------------------------------
# 3. Add the "Bay Count" sheet (after "Base Measurement").
# ---------------------------------------------------------------
$bayCount = $wb.Worksheets.Add($null, $baseMeasurement)
$bayCount.Name = "Bay Count"

$bayCount.Range("A1").Value2 = "KPI Name"
$bayCount.Range("A2").Value2 = "Warm Bays"
$bayCount.Range("A3").Value2 = "Cold Room Bays"
$bayCount.Range("A4").Value2 = "Cooler Door Count"

# formatting: header row bold/green, rows 2-4 bordered like other KPI sheets
$headerStyleSrc.Copy() | Out-Null
$bayCount.Range("A1").PasteSpecial(-4122) | Out-Null

$borderStyleSrc.Copy() | Out-Null
$bayCount.Range("A2:A4").PasteSpecial(-4122) | Out-Null

$bayCount.Columns.Item(1).ColumnWidth = 16.2112010796222

$bayCount.Rows.Item(1).RowHeight = 14.95
$bayCount.Rows.Item(2).RowHeight = 13.8
$bayCount.Rows.Item(3).RowHeight = 13.8
$bayCount.Rows.Item(4).RowHeight = 13.8

[void]$bayCount.Range("B2").Select()

# Bay Count becomes the active / selected sheet (tab) of the workbook.
$bayCount.Activate()

$excel.CutCopyMode = $false
